$wb = $excel.ActiveWorkbook

$seSheet = $wb.Worksheets.Item("SE")
$a3ssSheet = $wb.Worksheets.Item("A3SS")

# Data rows to append to SE sheet (sheet1), starting at row 2.
# Each inner array is one data row, column A..L.
$seRows = @(
    @("chr10", "TPD52L1", 31343011, 31343026, 0.0375962913753642, -0.2443458799799822, 0.2443458799799822, 0, 22.72291993720565, 9.4375, 29.18602825745682, "chr10:TPD52L1:31343011-31343026"),
    @("chr10", "TRMT11", 30591020, 30591102, 0.6811682297322682, 0.0459023371929579, 0.0459023371929579, 70.74698795180723, 92.98034076015728, 73.4578313253012, 116.7487986020096, "chr10:TRMT11:30591020-30591102"),
    @("chr10", "TPD52L1", 31341139, 31341178, 0.759639026899154, -0.0710106975581275, 0.0710106975581275, 3.425, 23.0496, 7.425, 29.6296, "chr10:TPD52L1:31341139-31341178"),
    @("chr10", "TPD52L1", 31338184, 31338245, 0.8905158238694521, 0.0510201852206924, 0.0510201852206924, 11.20967741935484, 23.00814332247557, 11.40322580645161, 29.82654723127036, "chr10:TPD52L1:31338184-31338245"),
    @("chr10", "ZUFSP", 33935088, 33935257, 0.9533697673964592, -0.0213551662639566, 0.0213551662639566, 41.96470588235294, 54.60894308943089, 50.81764705882353, 60.65121951219512, "chr10:ZUFSP:33935088-33935257")
)

# Write column-by-column (matches how the original file's shared-string
# table was built by the authoring tool).
for ($col = 1; $col -le 12; $col++) {
    for ($i = 0; $i -lt $seRows.Count; $i++) {
        $rowIdx = $i + 2
        $seSheet.Cells.Item($rowIdx, $col).Value = $seRows[$i][$col - 1]
    }
}

# Data row to append to A3SS sheet (sheet4), starting at row 2.
$a3ssRows = @(
    ,@("chr10", "NCOA7", 30691755, 30691787, 0.8698350074869731, 0.026160299348494, 0.026160299348494, 88.33333333333333, 95.53465184700885, 102.7272727272727, 123.4163125204315, "chr10:NCOA7:30691755-30691787")
)

for ($col = 1; $col -le 12; $col++) {
    for ($i = 0; $i -lt $a3ssRows.Count; $i++) {
        $rowIdx = $i + 2
        $a3ssSheet.Cells.Item($rowIdx, $col).Value = $a3ssRows[$i][$col - 1]
    }
}
